$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that sits right under the
#    H1 title (it duplicated content that is being moved to the very
#    end of the document, just above the image-prompt paragraph).
# ------------------------------------------------------------------
$metaIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Meta description*") {
        $metaIdx = $i
        break
    }
}
if ($metaIdx -gt 0) {
    $d.Paragraphs.Item($metaIdx).Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play BountyPop Free: Exciting
#    PopWins and Bonus Features" right before the final (italic image
#    prompt) paragraph.
# ------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIdx)
$startPos = $lastPara.Range.Start
$insertRange = $d.Range($startPos, $startPos)
$xmlNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newParaXml = "<w:p $xmlNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play BountyPop Free: Exciting PopWins and Bonus Features</w:t></w:r></w:p><w:p $xmlNs/>"
$null = $insertRange.InsertXML($newParaXml)

# InsertXML above splits the original last paragraph into two pieces,
# leaving a spurious empty paragraph between the new one and the
# original final paragraph -- remove it (its Range.Text is just the
# paragraph mark, not a zero-length string).
$emptyIdx = $lastIdx + 1
$emptyText = $d.Paragraphs.Item($emptyIdx).Range.Text
if ($emptyText.Length -le 1) {
    $d.Paragraphs.Item($emptyIdx).Range.Delete()
}

# ------------------------------------------------------------------
# 3) Swap the image-generation prompt text in the final paragraph for
#    the meta-description copy (minus the "Meta description: " label).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a cartoon-style feature image for BountyPop that features a happy Maya warrior with glasses. The image should be vibrant and eye-catching, featuring the Maya warrior surrounded by explosive gems and treasure chests, highlighting the adventurous pirate theme of the game. The image should encourage players to embrace the spirit of adventure and excitement while playing the game. The Maya warrior should be depicted as confident and carefree, perfectly embodying the attitude of players who are enjoying the game.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of BountyPop and play for free. Experience exciting PopWins feature, high volatility, and bonus games like Multiplier Wheel and Wheel of Bets.",
    2
)
